# Apply the "cryptos list" update produced by the scheduled GitHub Actions refresh
# (Sat Nov 11 13:56:00 UTC 2023): refreshed prices / 1h volume percentages, and
# swapped the ARBITRUM / FTXToken rows (43 <-> 44) to reflect the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.156.77'
$ws.Range("E2").Value = '  +0.03%  '

# Row 3
$ws.Range("D3").Value = '2.061.23'
$ws.Range("E3").Value = '  -1.38%  '

# Row 4
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").Value = '''250.45'
$ws.Range("E5").Value = '  -0.28%  '

# Row 6
$ws.Range("D6").Value = '''0.677'
$ws.Range("E6").Value = '  +3.22%  '

# Row 7
$ws.Range("D7").Value = '''59.35'
$ws.Range("E7").Value = '  +17.37%  '

# Row 8
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
$ws.Range("D9").Value = '''61.01'
$ws.Range("E9").Value = '  +0.33%  '

# Row 10
$ws.Range("E10").Value = '  +2.38%  '

# Row 11
$ws.Range("E11").Value = '  +7.57%  '

# Row 12
$ws.Range("E12").Value = '  +1.52%  '

# Row 13
$ws.Range("D13").Value = '''15.23'
$ws.Range("E13").Value = '  +0.50%  '

# Row 14
$ws.Range("D14").Value = '2.362.50'
$ws.Range("E14").Value = '  -1.36%  '

# Row 15
$ws.Range("D15").Value = '''0.818'
$ws.Range("E15").Value = '  -1.59%  '

# Row 16
$ws.Range("D16").Value = '''5.34'
$ws.Range("E16").Value = '  +4.70%  '

# Row 17
$ws.Range("D17").Value = '2.070.95'
$ws.Range("E17").Value = '  -1.04%  '

# Row 18
$ws.Range("D18").Value = '37.118.90'
$ws.Range("E18").Value = '  +0.09%  '

# Row 19
$ws.Range("D19").Value = '''75.31'
$ws.Range("E19").Value = '  +4.15%  '

# Row 20
$ws.Range("E20").Value = '  +12.26%  '

# Row 21
$ws.Range("D21").Value = '''14.51'
$ws.Range("E21").Value = '  +9.41%  '

# Row 22
$ws.Range("D22").Value = '''5.38'
$ws.Range("E22").Value = '  +3.34%  '

# Row 23
$ws.Range("D23").Value = '''239.29'
$ws.Range("E23").Value = '  -0.31%  '

# Row 24
$ws.Range("E24").Value = '  -0.02%  '

# Row 25
$ws.Range("D25").Value = '''2.45'
$ws.Range("E25").Value = '  -0.93%  '

# Row 26
$ws.Range("D26").Value = '''171.64'
$ws.Range("E26").Value = '  +1.44%  '

# Row 27
$ws.Range("E27").Value = '  -0.89%  '

# Row 28
$ws.Range("D28").Value = '''20.28'
$ws.Range("E28").Value = '  -3.98%  '

# Row 29
$ws.Range("E29").Value = '  +0.48%  '

# Row 30
$ws.Range("E30").Value = '  +2.63%  '

# Row 31
$ws.Range("E31").Value = '  +3.18%  '

# Row 32
$ws.Range("E32").Value = '  -4.27%  '

# Row 33
$ws.Range("D33").Value = '''0.0633'
$ws.Range("E33").Value = '  +4.47%  '

# Row 34
$ws.Range("E34").Value = '  +8.31%  '

# Row 35
$ws.Range("D35").Value = '''0.0888'
$ws.Range("E35").Value = '  -4.53%  '

# Row 36
$ws.Range("E36").Value = '  -0.10%  '

# Row 37
$ws.Range("E37").Value = '  +0.17%  '

# Row 38
$ws.Range("E38").Value = '  -3.57%  '

# Row 39
$ws.Range("D39").Value = '''0.112'
$ws.Range("E39").Value = '  +28.55%  '

# Row 40
$ws.Range("D40").Value = '''1.35'
$ws.Range("E40").Value = '  +2.56%  '

# Row 41
$ws.Range("D41").Value = '''18.44'
$ws.Range("E41").Value = '  +4.60%  '

# Row 42
$ws.Range("E42").Value = '  +0.80%  '

# Row 43 (ARBITRUM -> FTXToken)
$ws.Range("B43").Value = 'FTXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D43").Value = '''4.51'
$ws.Range("E43").Value = '  +32.45%  '

# Row 44 (FTXToken -> ARBITRUM)
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '''1.15'
$ws.Range("E44").Value = '  +0.14%  '

# Row 45
$ws.Range("D45").Value = '''97.48'
$ws.Range("E45").Value = '  -0.06%  '

# Row 46
$ws.Range("D46").Value = '''2.80'
$ws.Range("E46").Value = '  -0.30%  '

# Row 47
$ws.Range("D47").Value = '''4.49'
$ws.Range("E47").Value = '  +13.51%  '

# Row 48
$ws.Range("E48").Value = '  +11.54%  '

# Row 49
$ws.Range("D49").Value = '1.305.95'
$ws.Range("E49").Value = '  -0.19%  '

# Row 50
$ws.Range("E50").Value = '  -2.36%  '

# Row 51
$ws.Range("D51").Value = '''6.90'
$ws.Range("E51").Value = '  +0.09%  '

# Clear the implicit "quote prefix" style that Excel applies to the numeric-looking
# text cells above, so the cells keep the workbook's default (unstyled) formatting.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").Style = "Normal"

